$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2878.2
$ws.Range("I62").Value = 2540.2
$ws.Range("J62").Value = 3554.2
$ws.Range("K62").Value = 2540.2
$ws.Range("L62").Value = 3554.2
$ws.Range("M62").Value = -1916.2
$ws.Range("N62").Value = -4802.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2878.2
$ws.Range("I65").Value = 2540.2
$ws.Range("J65").Value = 3554.2
$ws.Range("K65").Value = 12701
$ws.Range("L65").Value = 17771
$ws.Range("M65").Value = -9581
$ws.Range("N65").Value = -24011

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3891.5715
$ws.Range("I74").Value = 3260.3
$ws.Range("J74").Value = 4465.4546
$ws.Range("K74").Value = 3260.3
$ws.Range("L74").Value = 4465.4546
$ws.Range("M74").Value = -2324.3
$ws.Range("N74").Value = -6337.4546

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3891.5715
$ws.Range("I77").Value = 3260.3
$ws.Range("J77").Value = 4465.4546
$ws.Range("K77").Value = 16301.5
$ws.Range("L77").Value = 22327.273
$ws.Range("M77").Value = -11621.5
$ws.Range("N77").Value = -31687.273

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 395.625
$ws.Range("I92").Value = 381.42856
$ws.Range("J92").Value = 495
$ws.Range("K92").Value = 381.42856
$ws.Range("L92").Value = 495
$ws.Range("M92").Value = 866.5714399999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1535.3846
$ws.Range("I98").Value = 1708.5714
$ws.Range("J98").Value = 1333.3334
$ws.Range("K98").Value = 1708.5714
$ws.Range("L98").Value = 1333.3334
$ws.Range("M98").Value = -210.5714

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1535.3846
$ws.Range("I122").Value = 1708.5714
$ws.Range("J122").Value = 1333.3334
$ws.Range("K122").Value = 5125.7142
$ws.Range("L122").Value = 4000.0002
$ws.Range("M122").Value = -2675.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2733.4
$ws.Range("I2").Value = 2860.1
$ws.Range("J2").Value = 2480
$ws.Range("K2").Value = 2860.1
$ws.Range("L2").Value = 2480
$ws.Range("M2").Value = -2747.1
$ws.Range("N2").Value = -2706

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3997.6711
$ws.Range("I32").Value = 2534.9275
$ws.Range("J32").Value = 18416.143
$ws.Range("K32").Value = 2534.9275
$ws.Range("L32").Value = 18416.143
$ws.Range("M32").Value = -2247.9275

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 5505
$ws.Range("I38").Value = 3010
$ws.Range("J38").Value = 8000
$ws.Range("K38").Value = 3010
$ws.Range("L38").Value = 8000
$ws.Range("M38").Value = -2543
$ws.Range("N38").Value = -8934

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2733.4
$ws.Range("I116").Value = 2860.1
$ws.Range("J116").Value = 2480
$ws.Range("K116").Value = 2860.1
$ws.Range("L116").Value = 2480
$ws.Range("M116").Value = -566.0999999999999
$ws.Range("N116").Value = -7068

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2926.8333
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2926.8333
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 8780.499899999999
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -13680.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2269.55
$ws.Range("I132").Value = 893.38464
$ws.Range("J132").Value = 4825.2856
$ws.Range("K132").Value = 2680.15392
$ws.Range("L132").Value = 14475.8568
$ws.Range("M132").Value = -150.1539199999997
$ws.Range("N132").Value = -19535.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2733.4
$ws.Range("I3").Value = 2860.1
$ws.Range("J3").Value = 2480
$ws.Range("K3").Value = 2860.1
$ws.Range("L3").Value = 2480
$ws.Range("M3").Value = -2746.1
$ws.Range("N3").Value = -2708

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 20000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 20000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 20000
$ws.Range("N33").Value = -20672

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2392.6924
$ws.Range("I86").Value = 2913.25
$ws.Range("J86").Value = 1559.8
$ws.Range("K86").Value = 2913.25
$ws.Range("L86").Value = 1559.8
$ws.Range("M86").Value = -1790.25
$ws.Range("N86").Value = -3805.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2392.6924
$ws.Range("I89").Value = 2913.25
$ws.Range("J89").Value = 1559.8
$ws.Range("K89").Value = 14566.25
$ws.Range("L89").Value = 7799
$ws.Range("M89").Value = -8950.25
$ws.Range("N89").Value = -19031

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4750
$ws.Range("I105").Value = 4333.3335
$ws.Range("J105").Value = 6000
$ws.Range("K105").Value = 4333.3335
$ws.Range("L105").Value = 6000
$ws.Range("M105").Value = -2586.3335
$ws.Range("N105").Value = -9494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 9262.556
$ws.Range("I35").Value = 1666.875
$ws.Range("J35").Value = 70028
$ws.Range("K35").Value = 1666.875
$ws.Range("L35").Value = 70028
$ws.Range("M35").Value = -1372.875
$ws.Range("N35").Value = -70616

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 2614.7058
$ws.Range("I133").Value = 3306.25
$ws.Range("J133").Value = 2000
$ws.Range("K133").Value = 9918.75
$ws.Range("L133").Value = 6000
$ws.Range("M133").Value = -4858.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2032.3334
$ws.Range("I31").Value = 531
$ws.Range("J31").Value = 5035
$ws.Range("K31").Value = 531
$ws.Range("L31").Value = 5035
$ws.Range("M31").Value = -239
$ws.Range("N31").Value = -5619

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 2032.3334
$ws.Range("I37").Value = 531
$ws.Range("J37").Value = 5035
$ws.Range("K37").Value = 531
$ws.Range("L37").Value = 5035
$ws.Range("M37").Value = -254
$ws.Range("N37").Value = -5589

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 5429.357
$ws.Range("I46").Value = 2520.5
$ws.Range("J46").Value = 5914.1665
$ws.Range("K46").Value = 2520.5
$ws.Range("L46").Value = 5914.1665
$ws.Range("M46").Value = -2364.5
$ws.Range("N46").Value = -6226.1665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2318.0488
$ws.Range("I122").Value = 1999.742
$ws.Range("J122").Value = 3304.8
$ws.Range("K122").Value = 5999.226
$ws.Range("L122").Value = 9914.400000000001
$ws.Range("M122").Value = -3549.226
$ws.Range("N122").Value = -14814.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4645.722
$ws.Range("I132").Value = 3541.1
$ws.Range("J132").Value = 6026.5
$ws.Range("K132").Value = 10623.3
$ws.Range("L132").Value = 18079.5
$ws.Range("M132").Value = -8093.299999999999
$ws.Range("N132").Value = -23139.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4632916.5
$ws.Range("I122").Value = 11112234
$ws.Range("J122").Value = 4832.7856
$ws.Range("K122").Value = 33336702
$ws.Range("L122").Value = 14498.3568
$ws.Range("M122").Value = -33334252
$ws.Range("N122").Value = -19398.3568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 15624.263
$ws.Range("I132").Value = 32422.572
$ws.Range("J132").Value = 5825.25
$ws.Range("K132").Value = 97267.716
$ws.Range("L132").Value = 17475.75
$ws.Range("M132").Value = -94737.716
$ws.Range("N132").Value = -22535.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 33339306
$ws.Range("I136").Value = 5617.1665
$ws.Range("J136").Value = 83339840
$ws.Range("K136").Value = 16851.4995
$ws.Range("L136").Value = 250019520
$ws.Range("M136").Value = -14301.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 251625
$ws.Range("I126").Value = 333833.34
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 1001500.02
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -999030.02

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2807.9167
$ws.Range("I132").Value = 2360.389
$ws.Range("J132").Value = 4150.5
$ws.Range("K132").Value = 7081.167
$ws.Range("L132").Value = 12451.5
$ws.Range("M132").Value = -4551.167
$ws.Range("N132").Value = -17511.5
